# Auto-generated Excel COM-interop script applying the profit-recalc edits
# described by the commit diff ("chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 62.42857
$ws.Range("I11").Value = 62.42857
$ws.Range("K11").Value = 62.42857
$ws.Range("M11").Value = 77.57142999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3440.3076
$ws.Range("I29").Value = 1230
$ws.Range("J29").Value = 4422.6665
$ws.Range("K29").Value = 3690
$ws.Range("L29").Value = 13267.9995
$ws.Range("M29").Value = -3409
$ws.Range("N29").Value = -13829.9995

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5062.913
$ws.Range("J40").Value = 7374.625
$ws.Range("L40").Value = 7374.625
$ws.Range("N40").Value = -7724.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6848.75
$ws.Range("J64").Value = 9180
$ws.Range("L64").Value = 9180
$ws.Range("N64").Value = -9676

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 6848.75
$ws.Range("J67").Value = 9180
$ws.Range("L67").Value = 9180
$ws.Range("N67").Value = -10896

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 615.375
$ws.Range("I92").Value = 703.8333
$ws.Range("K92").Value = 703.8333
$ws.Range("M92").Value = 544.1667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 5615.5
$ws.Range("I113").Value = 2981
$ws.Range("K113").Value = 2981
$ws.Range("M113").Value = 273

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3480.5625
$ws.Range("I137").Value = 2338.4
$ws.Range("K137").Value = 7015.200000000001
$ws.Range("M137").Value = -4465.200000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2665.5334
$ws.Range("I138").Value = 2498.7856
$ws.Range("J138").Value = 5000
$ws.Range("K138").Value = 7496.3568
$ws.Range("L138").Value = 15000
$ws.Range("M138").Value = -2356.3568
$ws.Range("N138").Value = -25280

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1736.1111
$ws.Range("I141").Value = 920
$ws.Range("J141").Value = 2756.25
$ws.Range("K141").Value = 2760
$ws.Range("L141").Value = 8268.75
$ws.Range("M141").Value = 2420
$ws.Range("N141").Value = -18628.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2501762
$ws.Range("J32").Value = 7146893
$ws.Range("L32").Value = 7146893
$ws.Range("N32").Value = -7147467

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 5057610.5
$ws.Range("J43").Value = 1743480.6
$ws.Range("L43").Value = 1743480.6
$ws.Range("N43").Value = -1744106.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 10049
$ws.Range("J54").Value = 10049
$ws.Range("L54").Value = 10049
$ws.Range("N54").Value = -11587

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3559.6875
$ws.Range("I61").Value = 2177.818
$ws.Range("K61").Value = 2177.818
$ws.Range("M61").Value = -1965.818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 735.4
$ws.Range("I110").Value = 706.1111
$ws.Range("K110").Value = 706.1111
$ws.Range("M110").Value = 1338.8889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3215.7778
$ws.Range("I122").Value = 3087.75
$ws.Range("K122").Value = 9263.25
$ws.Range("M122").Value = -6813.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3163.4211
$ws.Range("I132").Value = 3163.4211
$ws.Range("K132").Value = 9490.263300000001
$ws.Range("M132").Value = -6960.263300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3559.6875
$ws.Range("I136").Value = 2177.818
$ws.Range("K136").Value = 6533.454000000001
$ws.Range("M136").Value = -3983.454000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 28553.5
$ws.Range("I75").Value = 16404.666
$ws.Range("K75").Value = 16404.666
$ws.Range("M75").Value = -15468.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H78").Value = 28553.5
$ws.Range("I78").Value = 16404.666
$ws.Range("K78").Value = 49213.99800000001
$ws.Range("M78").Value = -44533.99800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5361.52
$ws.Range("I86").Value = 3990.375
$ws.Range("K86").Value = 3990.375
$ws.Range("M86").Value = -2867.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 5361.52
$ws.Range("I89").Value = 3990.375
$ws.Range("K89").Value = 19951.875
$ws.Range("M89").Value = -14335.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 38987
$ws.Range("I96").Value = 38987
$ws.Range("K96").Value = 38987
$ws.Range("M96").Value = -36241

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1558.75
$ws.Range("I134").Value = 1558.75
$ws.Range("K134").Value = 4676.25
$ws.Range("M134").Value = -2141.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 12461.667
$ws.Range("J141").Value = 9923.333000000001
$ws.Range("L141").Value = 9923.333000000001
$ws.Range("N141").Value = -20283.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5829.433
$ws.Range("I31").Value = 1328.8182
$ws.Range("K31").Value = 1328.8182
$ws.Range("M31").Value = -1033.8182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5829.433
$ws.Range("I34").Value = 1328.8182
$ws.Range("K34").Value = 1328.8182
$ws.Range("M34").Value = -1126.8182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 24750
$ws.Range("J57").Value = 24750
$ws.Range("L57").Value = 24750
$ws.Range("N57").Value = -25870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 829
$ws.Range("I134").Value = 829
$ws.Range("K134").Value = 2487
$ws.Range("M134").Value = 48

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 629.8
$ws.Range("J92").Value = 699.6667
$ws.Range("L92").Value = 2099.0001
$ws.Range("N92").Value = -4595.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 481.66666
$ws.Range("J97").Value = 497.5
$ws.Range("L97").Value = 1492.5
$ws.Range("N97").Value = -2484.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 2816.1667
$ws.Range("J117").Value = 2979.4
$ws.Range("L117").Value = 8938.200000000001
$ws.Range("N117").Value = -15822.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 5000
$ws.Range("J137").Value = 5000
$ws.Range("L137").Value = 15000
$ws.Range("N137").Value = -25200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2500
$ws.Range("I70").Value = 2500
$ws.Range("K70").Value = 2500
$ws.Range("M70").Value = -2230

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 2500
$ws.Range("I73").Value = 2500
$ws.Range("K73").Value = 2500
$ws.Range("M73").Value = -1564

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1700.8182
$ws.Range("I102").Value = 1790
$ws.Range("K102").Value = 1790
$ws.Range("M102").Value = -168

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1718.4286
$ws.Range("I122").Value = 1713.1666
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 5139.4998
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -2689.4998
$ws.Range("N122").Value = -10150

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3552.5862
$ws.Range("J132").Value = 5989.857
$ws.Range("L132").Value = 17969.571
$ws.Range("N132").Value = -23029.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 4000
$ws.Range("K122").Value = 12000
$ws.Range("M122").Value = -9550

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3397.25
$ws.Range("I132").Value = 3397.25
$ws.Range("K132").Value = 10191.75
$ws.Range("M132").Value = -7661.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 39999
$ws.Range("J69").Value = 39999
$ws.Range("L69").Value = 39999
$ws.Range("N69").Value = -41497

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 39999
$ws.Range("J72").Value = 39999
$ws.Range("L72").Value = 119997
$ws.Range("N72").Value = -127485

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 37000
$ws.Range("J112").Value = 37000
$ws.Range("L112").Value = 37000
$ws.Range("N112").Value = -39954

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5456.478
$ws.Range("I126").Value = 3608.5
$ws.Range("J126").Value = 7472.4546
$ws.Range("K126").Value = 10825.5
$ws.Range("L126").Value = 22417.3638
$ws.Range("M126").Value = -8355.5
$ws.Range("N126").Value = -27357.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2506.6
$ws.Range("I132").Value = 2455.75
$ws.Range("K132").Value = 7367.25
$ws.Range("M132").Value = -4837.25
